$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.281.44"
$ws.Range("E2").Value = "  +5.48%  "
$ws.Range("D3").Value = "1.919.08"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9996"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "254.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9996"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5184"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.36%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "46.12"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +7.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2994"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +7.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06893"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.98%  "
$ws.Range("D11").Value = "1.918.79"
$ws.Range("E11").Value = "  +6.09%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "17.55"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.45%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07355"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6906"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "88.13"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.939"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.91%  "
$ws.Range("D17").Value = "30.289.02"
$ws.Range("E17").Value = "  +5.60%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008250"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +12.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9984"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.13"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.04%  "
$ws.Range("D21").Value = "2.163.93"
$ws.Range("E21").Value = "  +6.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9995"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.877"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.57%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.790"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +8.96%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.223"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.73%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "147.16"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "139.93"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +25.77%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.35"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +8.41%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.017"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.379"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.308"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.97%  "
$ws.Range("E32").Value = "  +6.30%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.037"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.94%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05149"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.75%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.164"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.53%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7228"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.47%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.685"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.69%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.336"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.80%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.857"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.12%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9769"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.72%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01707"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.82%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.189"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4356"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.73%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "106.51"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.40%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9989"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.705"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.45%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1282"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.53%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05737"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.566"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.98%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.38"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.39%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3858"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.00%  "
